$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "What is the status of order number 1597?"
$ws.Range("C2").Value = "{'orders': [{'Product ID': 43558, 'ProductName': 'bosch kir21vf30g serie 4 built in single door fridges 88cm height white', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1597, 'CustomerID': 2244, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-14 12:24:37.765490'}]}"

$ws.Range("B3").Value = "Can I get an update on order 1728?"
$ws.Range("C3").Value = "{'orders': [{'Product ID': 16687, 'ProductName': 'hp intel xeon quad core e5506 2.13ghz fio kit 2.13ghz 4mb l2', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 1728, 'CustomerID': 2421, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-09-30 12:24:37.765490'}]}"

$ws.Range("B4").Value = "What's the latest on my order, 3421?"
$ws.Range("C4").Value = "{'orders': [{'Product ID': 22065, 'ProductName': 'sony dsc rx0 1.0 type sensor ultra compact camera with waterproof and shockproof design black', 'Category': 'Digital Cameras', 'Category ID': 2617, 'OrderID': 3421, 'CustomerID': 1607, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-27 12:24:37.765490'}]}"

$ws.Range("B5").Value = "Could you check the status of order 3117 for me, please?"
$ws.Range("C5").Value = "{'orders': [{'Product ID': 29707, 'ProductName': 'bosch exclusiv made in germany supersilence geschirrsp ler 45 cm integrierbar edelstahl spi66ts00d', 'Category': 'Dishwashers', 'Category ID': 2619, 'OrderID': 3117, 'CustomerID': 880, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-10-12 12:24:37.765490'}]}"

$ws.Range("B6").Value = "Is there any information available about order number 1479?"
$ws.Range("C6").Value = "{'orders': [{'Product ID': 26272, 'ProductName': 'de longhi standard microwave p80d20 stainless steel', 'Category': 'Microwaves', 'Category ID': 2618, 'OrderID': 1479, 'CustomerID': 1353, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-10-09 12:24:37.765490'}]}"

$ws.Range("B7").Value = "What is the status of order number 2431?"
$ws.Range("C7").Value = "{'orders': [{'Product ID': 2040, 'ProductName': 'xiaomi redmi 6a 32gb lte black', 'Category': 'Mobile Phones', 'Category ID': 2612, 'OrderID': 1694, 'CustomerID': 304, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-13 12:24:37.765490'}]}"

$ws.Range("B8").Value = "Can I get an update on order 2179?"
$ws.Range("C8").Value = "{'orders': [{'Product ID': 36702, 'ProductName': 'hotpoint hf1801efaa 54cm wide frost free integrated upright in column freezer white', 'Category': 'Freezers', 'Category ID': 2621, 'OrderID': 1417, 'CustomerID': 1958, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-09-29 12:24:37.765490'}]}"

$ws.Range("B9").Value = "What's the latest on my order, 2681?"
$ws.Range("C9").Value = "{'orders': [{'Product ID': 41945, 'ProductName': 'bosch einbau k hl /gefrier kombination smartcool flachscharnier technik profi t rd mpfung kis77ad40 eek a', 'Category': 'Fridge Freezers', 'Category ID': 2622, 'OrderID': 1326, 'CustomerID': 2427, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-27 12:24:37.765490'}]}"

$ws.Range("B10").Value = "Could you check the status of order 2130 for me, please?"
$ws.Range("C10").Value = "{'orders': [{'Product ID': 3359, 'ProductName': 'energizer energy 100 sim mobile phone double sim bluetooth 2g ip54 black', 'Category': 'Mobile Phones', 'Category ID': 2612, 'OrderID': 1202, 'CustomerID': 142, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-06 12:24:37.765490'}]}"

$ws.Range("B11").Value = "Is there any information available about order number 2400?"
$ws.Range("C11").Value = "{'orders': [{'Product ID': 15315, 'ProductName': 'intel core i3 3250 3.5 ghz 2 cores 4 threads 3 mb cache lga1155 socket oem', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 2194, 'CustomerID': 334, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-06 12:24:37.765490'}]}"

$wb.Save()
